$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new dates/strings for weeks 7, 8 and 9 (rows 6, 7 and 8)
$ws.Range("C6").Value = "16.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D6").Value = "18.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."

$ws.Range("C7").Value = "23.02: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D7").Value = "25.02: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

$ws.Range("C8").Value = "02.03: Kontakttime, kursansvarlig tilgjengelig på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09)."
$ws.Range("D8").Value = "04.03: Oppgaveseminar på [Zoom](https://nhh.zoom.us/j/63633653066?pwd=cTVNV0JvOXl4NnUrMHVKdkw2b0pCZz09). Se \@ref(seminar) for oppgaver."

# Update the active cell selection to D9
$ws.Range("D9").Select()
